$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values in C3 and D3: the "Pedro-Resistencia mecanica" class
# moves from Tuesday (C3) to Wednesday (D3) for the 7:50 slot.
$ws.Range("C3").Value = "-"
$ws.Range("D3").Value = "Pedro-Resistencia mecanica"
